$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three fresh weekly rows before the existing row 44 data block.
# The prior rows 44:48 shift down to 47:51 automatically.
$ws.Rows("44:46").Insert()

# New row 44 - Espárragos, Banquete, "Sin especificar", Provincia de Linares
$ws.Cells.Item(44, 1).Value = 9
$ws.Cells.Item(44, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(44, 3).Value = "Metropolitana"
$ws.Cells.Item(44, 4).Value = 44474
$ws.Cells.Item(44, 5).Value = 13
$ws.Cells.Item(44, 6).Value = 300000000
$ws.Cells.Item(44, 7).Value = "Espárragos"
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 9).Value = "Banquete"
$ws.Cells.Item(44, 10).Value = 196
$ws.Cells.Item(44, 11).Value = 1500
$ws.Cells.Item(44, 12).Value = 1500
$ws.Cells.Item(44, 13).Value = 1500
$ws.Cells.Item(44, 14).Value = "$/kilo"
$ws.Cells.Item(44, 15).Value = "Provincia de Linares"
$ws.Cells.Item(44, 16).Value = 1500
$ws.Cells.Item(44, 17).Value = 1
$ws.Cells.Item(44, 18).Value = "Hortaliza"

# New row 45 - Espárragos, Primera, "Sin especificar", Provincia de Linares
$ws.Cells.Item(45, 1).Value = 9
$ws.Cells.Item(45, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(45, 3).Value = "Metropolitana"
$ws.Cells.Item(45, 4).Value = 44474
$ws.Cells.Item(45, 5).Value = 13
$ws.Cells.Item(45, 6).Value = 300000000
$ws.Cells.Item(45, 7).Value = "Espárragos"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 340
$ws.Cells.Item(45, 11).Value = 1300
$ws.Cells.Item(45, 12).Value = 1300
$ws.Cells.Item(45, 13).Value = 1300
$ws.Cells.Item(45, 14).Value = "$/kilo"
$ws.Cells.Item(45, 15).Value = "Provincia de Linares"
$ws.Cells.Item(45, 16).Value = 1300
$ws.Cells.Item(45, 17).Value = 1
$ws.Cells.Item(45, 18).Value = "Hortaliza"

# New row 46 - Espárragos, Segunda, "Sin especificar", Provincia de Linares
$ws.Cells.Item(46, 1).Value = 9
$ws.Cells.Item(46, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(46, 3).Value = "Metropolitana"
$ws.Cells.Item(46, 4).Value = 44474
$ws.Cells.Item(46, 5).Value = 13
$ws.Cells.Item(46, 6).Value = 300000000
$ws.Cells.Item(46, 7).Value = "Espárragos"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Segunda"
$ws.Cells.Item(46, 10).Value = 160
$ws.Cells.Item(46, 11).Value = 1000
$ws.Cells.Item(46, 12).Value = 1000
$ws.Cells.Item(46, 13).Value = 1000
$ws.Cells.Item(46, 14).Value = "$/kilo"
$ws.Cells.Item(46, 15).Value = "Provincia de Linares"
$ws.Cells.Item(46, 16).Value = 1000
$ws.Cells.Item(46, 17).Value = 1
$ws.Cells.Item(46, 18).Value = "Hortaliza"
